# Auto-generated cell updates applying the cryptos list refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.418.03"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "3.511.61"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  +0.01%  "
$sStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.39"
$ws.Range("D5").Style = $sStyle
$ws.Range("E5").Value = "  +0.91%  "
$sStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.62"
$ws.Range("D6").Style = $sStyle
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("E9").Value = "  +6.10%  "
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("E11").Value = "  +3.68%  "
$ws.Range("D12").Value = "4.108.86"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("D15").Value = "3.511.29"
$ws.Range("E15").Value = "  +0.36%  "
$sStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.79"
$ws.Range("D16").Style = $sStyle
$ws.Range("E16").Value = "  +3.09%  "
$ws.Range("D17").Value = "64.402.71"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("E18").Value = "  -0.24%  "
$sStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.65"
$ws.Range("D19").Style = $sStyle
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("E20").Value = "  +2.12%  "
$sStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "396.24"
$ws.Range("D21").Style = $sStyle
$ws.Range("E21").Value = "  +2.96%  "
$ws.Range("E22").Value = "  +2.05%  "
$ws.Range("D23").Value = "3.651.18"
$ws.Range("E23").Value = "  +0.36%  "
$sStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.68"
$ws.Range("D24").Style = $sStyle
$ws.Range("E24").Value = "  +0.86%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("E27").Value = "  +3.63%  "
$ws.Range("E28").Value = "  +0.02%  "
$sStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.43"
$ws.Range("D29").Style = $sStyle
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("E30").Value = "  +1.25%  "
$sStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.25"
$ws.Range("D31").Style = $sStyle
$ws.Range("E31").Value = "  +0.35%  "
$ws.Range("E32").Value = "  -4.15%  "
$ws.Range("E33").Value = "  +6.61%  "
$ws.Range("D34").Value = "3.539.03"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("E36").Value = "  -0.62%  "
$sStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.36"
$ws.Range("D37").Style = $sStyle
$ws.Range("E37").Value = "  +1.42%  "
$ws.Range("E38").Value = "  +2.11%  "
$sStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.57"
$ws.Range("D39").Style = $sStyle
$ws.Range("E39").Value = "  +1.29%  "
$sStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "167.20"
$ws.Range("D40").Style = $sStyle
$ws.Range("E40").Value = "  +3.12%  "
$ws.Range("E41").Value = "  +0.99%  "
$sStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.811"
$ws.Range("D42").Style = $sStyle
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("E43").Value = "  +0.01%  "
$sStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.46"
$ws.Range("D44").Style = $sStyle
$ws.Range("E44").Value = "  +0.96%  "
$sStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.82"
$ws.Range("D45").Style = $sStyle
$ws.Range("E45").Value = "  -3.83%  "
$ws.Range("E46").Value = "  +1.24%  "
$sStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.18"
$ws.Range("D47").Style = $sStyle
$ws.Range("E47").Value = "  -2.77%  "
$sStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.82"
$ws.Range("D48").Style = $sStyle
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("D49").Value = "2.377.32"
$ws.Range("E49").Value = "  -3.80%  "
$sStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.909"
$ws.Range("D50").Style = $sStyle
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("E51").Value = "  +0.41%  "
